$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1235.6364
$ws.Range("I19").Value = 1137.25
$ws.Range("K19").Value = 1137.25
$ws.Range("M19").Value = -962.25
$ws.Range("H62").Value = 8676.117
$ws.Range("I62").Value = 8676.117
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 8676.117
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -8052.117
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 8676.117
$ws.Range("I65").Value = 8676.117
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 43380.585
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -40260.585
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 10589.25
$ws.Range("I69").Value = 8344.546
$ws.Range("J69").Value = 13332.777
$ws.Range("K69").Value = 25033.638
$ws.Range("L69").Value = 39998.331
$ws.Range("M69").Value = -24159.638
$ws.Range("N69").Value = -41746.331
$ws.Range("H72").Value = 10589.25
$ws.Range("I72").Value = 8344.546
$ws.Range("J72").Value = 13332.777
$ws.Range("K72").Value = 75100.914
$ws.Range("L72").Value = 119994.993
$ws.Range("M72").Value = -70732.914
$ws.Range("N72").Value = -128730.993
$ws.Range("H86").Value = 4956.3335
$ws.Range("I86").Value = 4997.8887
$ws.Range("K86").Value = 4997.8887
$ws.Range("M86").Value = -3874.8887
$ws.Range("H88").Value = 9227.416999999999
$ws.Range("J88").Value = 8581.111000000001
$ws.Range("L88").Value = 8581.111000000001
$ws.Range("N88").Value = -9393.111000000001
$ws.Range("H89").Value = 4956.3335
$ws.Range("I89").Value = 4997.8887
$ws.Range("K89").Value = 24989.4435
$ws.Range("M89").Value = -19373.4435
$ws.Range("H91").Value = 9227.416999999999
$ws.Range("J91").Value = 8581.111000000001
$ws.Range("L91").Value = 8581.111000000001
$ws.Range("N91").Value = -11389.111
$ws.Range("H131").Value = 3496.2856
$ws.Range("I131").Value = 3496.2856
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 10488.8568
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -5448.856800000001
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 10219.904
$ws.Range("I132").Value = 10169.294
$ws.Range("J132").Value = 10435
$ws.Range("K132").Value = 30507.882
$ws.Range("L132").Value = 31305
$ws.Range("M132").Value = -27977.882
$ws.Range("N132").Value = -36365
$ws.Range("H138").Value = 20888.924
$ws.Range("I138").Value = 24007.143
$ws.Range("K138").Value = 72021.429
$ws.Range("M138").Value = -66881.429

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15472.975
$ws.Range("I32").Value = 15640.598
$ws.Range("K32").Value = 15640.598
$ws.Range("M32").Value = -15353.598
$ws.Range("H43").Value = 18188
$ws.Range("J43").Value = 18188
$ws.Range("L43").Value = 18188
$ws.Range("N43").Value = -18814
$ws.Range("H46").Value = 30075.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 30075.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 30075.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -30713.5
$ws.Range("H132").Value = 28943.525
$ws.Range("I132").Value = 38035.93
$ws.Range("J132").Value = 3484.8
$ws.Range("K132").Value = 114107.79
$ws.Range("L132").Value = 10454.4
$ws.Range("M132").Value = -111577.79
$ws.Range("N132").Value = -15514.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4724.278
$ws.Range("I105").Value = 4504.9287
$ws.Range("K105").Value = 4504.9287
$ws.Range("M105").Value = -2757.9287

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 1887.6818
$ws.Range("I105").Value = 1834.1904
$ws.Range("K105").Value = 1834.1904
$ws.Range("M105").Value = -87.19039999999995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 545.0625
$ws.Range("I8").Value = 545.0625
$ws.Range("K8").Value = 1635.1875
$ws.Range("M8").Value = -1496.1875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3674.1428
$ws.Range("I80").Value = 3228.2222
$ws.Range("J80").Value = 4476.8
$ws.Range("K80").Value = 3228.2222
$ws.Range("L80").Value = 4476.8
$ws.Range("M80").Value = -2230.2222
$ws.Range("N80").Value = -6472.8
$ws.Range("H83").Value = 3674.1428
$ws.Range("I83").Value = 3228.2222
$ws.Range("J83").Value = 4476.8
$ws.Range("K83").Value = 16141.111
$ws.Range("L83").Value = 22384
$ws.Range("M83").Value = -11149.111
$ws.Range("N83").Value = -32368

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4020
$ws.Range("I7").Value = 4020
$ws.Range("K7").Value = 4020
$ws.Range("M7").Value = -3908
$ws.Range("H22").Value = 41583.105
$ws.Range("I22").Value = 93290.664
$ws.Range("J22").Value = 2802.4375
$ws.Range("K22").Value = 93290.664
$ws.Range("L22").Value = 2802.4375
$ws.Range("M22").Value = -92995.664
$ws.Range("N22").Value = -3392.4375
$ws.Range("H27").Value = 41583.105
$ws.Range("I27").Value = 93290.664
$ws.Range("J27").Value = 2802.4375
$ws.Range("K27").Value = 93290.664
$ws.Range("L27").Value = 2802.4375
$ws.Range("M27").Value = -93183.664
$ws.Range("N27").Value = -3016.4375
$ws.Range("H82").Value = 2671.4285
$ws.Range("I82").Value = 2566.7778
$ws.Range("J82").Value = 2859.8
$ws.Range("K82").Value = 2566.7778
$ws.Range("L82").Value = 2859.8
$ws.Range("M82").Value = -2205.7778
$ws.Range("N82").Value = -3581.8
$ws.Range("H85").Value = 2671.4285
$ws.Range("I85").Value = 2566.7778
$ws.Range("J85").Value = 2859.8
$ws.Range("K85").Value = 2566.7778
$ws.Range("L85").Value = 2859.8
$ws.Range("M85").Value = -1318.7778
$ws.Range("N85").Value = -5355.8
$ws.Range("H126").Value = 4020
$ws.Range("I126").Value = 4020
$ws.Range("K126").Value = 12060
$ws.Range("M126").Value = -9590

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2837.2856
$ws.Range("I81").Value = 2928.1667
$ws.Range("J81").Value = 2292
$ws.Range("K81").Value = 5856.3334
$ws.Range("L81").Value = 4584
$ws.Range("M81").Value = -4795.3334
$ws.Range("N81").Value = -6706
$ws.Range("H84").Value = 2837.2856
$ws.Range("I84").Value = 2928.1667
$ws.Range("J84").Value = 2292
$ws.Range("K84").Value = 29281.667
$ws.Range("L84").Value = 22920
$ws.Range("M84").Value = -23977.667
$ws.Range("N84").Value = -33528
$ws.Range("H96").Value = 3427.818
$ws.Range("I96").Value = 2474.8
$ws.Range("K96").Value = 2474.8
$ws.Range("M96").Value = -1101.8
